# Insert a new weekly price record for "Haba" (Femacal de La Calera) above
# the current row 72, shifting the existing rows 72:87 down to 73:88.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

$newRow = 72
$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item($newRow, 3).Value = 'Coquimbo'
$ws.Cells.Item($newRow, 4).Value = 44505
$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = 100112026
$ws.Cells.Item($newRow, 7).Value = 'Haba'
$ws.Cells.Item($newRow, 8).Value = 'Sin especificar'
$ws.Cells.Item($newRow, 9).Value = 'Primera'
$ws.Cells.Item($newRow, 10).Value = 75
$ws.Cells.Item($newRow, 11).Value = 7000
$ws.Cells.Item($newRow, 12).Value = 8000
$ws.Cells.Item($newRow, 13).Value = 7467
$ws.Cells.Item($newRow, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($newRow, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item($newRow, 16).Value = 299
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = 'Hortaliza'
